$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row after the current last data row (35) so the table
#     grows from 20 periods to 21 periods. The footer rows (old 40/41)
#     shift down automatically to 41/42, including their merged cells.
$ws.Rows.Item(36).Insert()

# Row 36 should inherit the special "last row" bottom-border style that
# row 35 currently has, so copy row 35's full formatting+values there first.
$ws.Range("B35:J35").Copy($ws.Range("B36:J36"))

# Row 35 now becomes a regular data row, so restyle it like row 34 (a
# normal interior row) before we overwrite its values below.
$ws.Range("B34:J34").Copy($ws.Range("B35:J35"))

# --- Update the summary figures: one more overdue period (21) and the
#     additional Valor Mora for the newly added period (722448 + 52000).
$ws.Range("E11").Value = 774448
$ws.Range("F13").Value = 21

# --- Rewrite the detail table (rows 16-36) in its new sorted order,
#     grouped per worker, with the extra trailing period (2508) appended
#     for DAIRO RHENALS VALERO.
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1050969465"
$ws.Range("D16").Value = "JUAN CARLOS RAMOS PRENS"
$ws.Range("E16").Value = "1707"
$ws.Range("F16").Value = 14755
$ws.Range("G16").Value = 781242

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1143389736"
$ws.Range("D17").Value = "DANIELA ESTHER MUÑOZ BARRIOS"
$ws.Range("E17").Value = "1710"
$ws.Range("F17").Value = 14755
$ws.Range("G17").Value = 1300000

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1143389736"
$ws.Range("D18").Value = "DANIELA ESTHER MUÑOZ BARRIOS"
$ws.Range("E18").Value = "1711"
$ws.Range("F18").Value = 29509
$ws.Range("G18").Value = 1300000

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1143389736"
$ws.Range("D19").Value = "DANIELA ESTHER MUÑOZ BARRIOS"
$ws.Range("E19").Value = "1712"
$ws.Range("F19").Value = 29509
$ws.Range("G19").Value = 1300000

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1143389736"
$ws.Range("D20").Value = "DANIELA ESTHER MUÑOZ BARRIOS"
$ws.Range("E20").Value = "1801"
$ws.Range("F20").Value = 29509
$ws.Range("G20").Value = 1300000

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "1143389736"
$ws.Range("D21").Value = "DANIELA ESTHER MUÑOZ BARRIOS"
$ws.Range("E21").Value = "1802"
$ws.Range("F21").Value = 29509
$ws.Range("G21").Value = 1300000

$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "1143389736"
$ws.Range("D22").Value = "DANIELA ESTHER MUÑOZ BARRIOS"
$ws.Range("E22").Value = "1803"
$ws.Range("F22").Value = 29509
$ws.Range("G22").Value = 1300000

$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "1143389736"
$ws.Range("D23").Value = "DANIELA ESTHER MUÑOZ BARRIOS"
$ws.Range("E23").Value = "1804"
$ws.Range("F23").Value = 29509
$ws.Range("G23").Value = 1300000

$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = "1143389736"
$ws.Range("D24").Value = "DANIELA ESTHER MUÑOZ BARRIOS"
$ws.Range("E24").Value = "1805"
$ws.Range("F24").Value = 29509
$ws.Range("G24").Value = 1300000

$ws.Range("B25").Value = "CC"
$ws.Range("C25").Value = "1051891826"
$ws.Range("D25").Value = "YAID LUZ DIAZ TORRES"
$ws.Range("E25").Value = "1809"
$ws.Range("F25").Value = 1042
$ws.Range("G25").Value = 781242

$ws.Range("B26").Value = "CC"
$ws.Range("C26").Value = "9147800"
$ws.Range("D26").Value = "DAIRO RHENALS VALERO"
$ws.Range("E26").Value = "2410"
$ws.Range("F26").Value = 17333
$ws.Range("G26").Value = 1300000

$ws.Range("B27").Value = "CC"
$ws.Range("C27").Value = "9147800"
$ws.Range("D27").Value = "DAIRO RHENALS VALERO"
$ws.Range("E27").Value = "2411"
$ws.Range("F27").Value = 52000
$ws.Range("G27").Value = 1300000

$ws.Range("B28").Value = "CC"
$ws.Range("C28").Value = "9147800"
$ws.Range("D28").Value = "DAIRO RHENALS VALERO"
$ws.Range("E28").Value = "2412"
$ws.Range("F28").Value = 52000
$ws.Range("G28").Value = 1300000

$ws.Range("B29").Value = "CC"
$ws.Range("C29").Value = "9147800"
$ws.Range("D29").Value = "DAIRO RHENALS VALERO"
$ws.Range("E29").Value = "2501"
$ws.Range("F29").Value = 52000
$ws.Range("G29").Value = 1300000

$ws.Range("B30").Value = "CC"
$ws.Range("C30").Value = "9147800"
$ws.Range("D30").Value = "DAIRO RHENALS VALERO"
$ws.Range("E30").Value = "2502"
$ws.Range("F30").Value = 52000
$ws.Range("G30").Value = 1300000

$ws.Range("B31").Value = "CC"
$ws.Range("C31").Value = "9147800"
$ws.Range("D31").Value = "DAIRO RHENALS VALERO"
$ws.Range("E31").Value = "2503"
$ws.Range("F31").Value = 52000
$ws.Range("G31").Value = 1300000

$ws.Range("B32").Value = "CC"
$ws.Range("C32").Value = "9147800"
$ws.Range("D32").Value = "DAIRO RHENALS VALERO"
$ws.Range("E32").Value = "2504"
$ws.Range("F32").Value = 52000
$ws.Range("G32").Value = 1300000

$ws.Range("B33").Value = "CC"
$ws.Range("C33").Value = "9147800"
$ws.Range("D33").Value = "DAIRO RHENALS VALERO"
$ws.Range("E33").Value = "2505"
$ws.Range("F33").Value = 52000
$ws.Range("G33").Value = 1300000

$ws.Range("B34").Value = "CC"
$ws.Range("C34").Value = "9147800"
$ws.Range("D34").Value = "DAIRO RHENALS VALERO"
$ws.Range("E34").Value = "2506"
$ws.Range("F34").Value = 52000
$ws.Range("G34").Value = 1300000

$ws.Range("B35").Value = "CC"
$ws.Range("C35").Value = "9147800"
$ws.Range("D35").Value = "DAIRO RHENALS VALERO"
$ws.Range("E35").Value = "2507"
$ws.Range("F35").Value = 52000
$ws.Range("G35").Value = 1300000

$ws.Range("B36").Value = "CC"
$ws.Range("C36").Value = "9147800"
$ws.Range("D36").Value = "DAIRO RHENALS VALERO"
$ws.Range("E36").Value = "2508"
$ws.Range("F36").Value = 52000
$ws.Range("G36").Value = 1300000
